$wb = $excel.ActiveWorkbook

# --- BOM sheet updates ---
$bom = $wb.Worksheets.Item("BOM")
$bom.Range("B20").Value = 24.68
$bom.Range("A30").Value = "1/2`" x 19`" aluminum sqaure bar"
$bom.Range("B30").Value = 7.95

# --- Insert a new "Extrusions" sheet between "BOM" and "Bulk Hardware" ---
$extrusions = $wb.Worksheets.Add($null, $bom)
$extrusions.Name = "Extrusions"

$extrusions.Range("A1").Value = "Extrusions"

$extrusions.Range("A2").Value = "Size"
$extrusions.Range("B2").Value = "Length"
$extrusions.Range("C2").Value = "Quantity"

$extrusions.Range("A3").Value = "4080 C-Channel"
$extrusions.Range("B3").Value = "914mm"
$extrusions.Range("C3").Value = 2

$extrusions.Range("A4").Value = "2040 V-slot"
$extrusions.Range("B4").Value = "520mm"
$extrusions.Range("C4").Value = 3

$extrusions.Range("A5").Value = "2040 V-slot"
$extrusions.Range("B5").Value = "700mm"
$extrusions.Range("C5").Value = 2

$extrusions.Range("A6").Value = "2040 V-slot"
$extrusions.Range("B6").Value = "460mm"
$extrusions.Range("C6").Value = 1

$extrusions.Range("A7").Value = "2040 V-slot"
$extrusions.Range("B7").Value = "453mm"
$extrusions.Range("C7").Value = 1

$extrusions.Range("A8").Value = "2020 V-slot"
$extrusions.Range("B8").Value = "315mm"
$extrusions.Range("C8").Value = 2

$extrusions.Range("A1").ColumnWidth = 13.83

# --- sheet view / tab selection housekeeping ---
$pinSheet = $wb.Worksheets.Item("34pin connector")
$pinSheet.Select()
$bom.Select()
$bom.Range("B22").Select()

Write-Host "Done"
